$wb = $excel.ActiveWorkbook

# --- NewImportLogic sheet ---
$ws1 = $wb.Worksheets.Item("NewImportLogic")
$ws1.Activate()
$ws1.Range("H4").Value = "StandardExcelReport-QOL_and_ECON - UtilityOutcome-Economic-2023_"
$ws1.Range("H4").Select() | Out-Null

# --- OldImportLogic sheet ---
$ws2 = $wb.Worksheets.Item("OldImportLogic")
$ws2.Activate()
$ws2.Range("H2").Value = "StandardExcelReport-Alkermes - Melanoma-Economic-2023_"
$ws2.Range("G8").Select() | Out-Null
